# Update the first JSON payload cell: quote the userId value ("5" instead of 5)
# and drop the stray trailing space before the closing brace.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "{\`"title\`": \`"jonnyalexfoo1\`", \`"body\`": \`"bar\`", \`"userId\`": \`"5\`"}"

# The second sample row (row 2) is no longer needed - remove it entirely.
$ws.Rows.Item(2).Delete()

# Column C was manually touched, so Excel no longer considers its width to be
# an auto "best fit" value - nudge the width back to (effectively) the same
# size while clearing that best-fit flag.
$ws.Columns.Item(3).ColumnWidth = 55
